# إضافة حدث جديد في Card11 by HOSSAM at 2025-12-08 11:41:30
# Adds a new service-log row (row 14) to the "Card11" sheet describing a
# coiler timing-belt (1270) replacement event.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card11")

$row = 14

# Column A ("card") keeps the same machine/card number as every other row,
# stored as text (like the rest of the column) rather than a number.
$ws.Cells.Item($row, 1).Value = "'11"
$ws.Cells.Item($row, 1).Style = "Normal"

# Columns B-K (tonnage + wear-part checklist) are left blank for this entry.
$ws.Cells.Item($row, 12).Value = "24/3/2025"
$ws.Cells.Item($row, 13).Value = "قطع سير كويلر مسنن 1270"
$ws.Cells.Item($row, 14).Value = "تم تغير سير 1270"
$ws.Cells.Item($row, 15).Value = "فني"
